$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from row 135 into the 5 new rows 136-140
$ws.Range("A135:AC135").Copy() | Out-Null
$ws.Range("A136:AC140").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 136
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 6787896
$ws.Range("C136").Value = "Croatia HNL"
$ws.Range("D136").Value = "Croatia HNL"
$ws.Range("E136").Value = 45380.49305555555
$ws.Range("F136").Value = "Istra 1961"
$ws.Range("G136").Value = "NK Lokomotiva Zagreb"
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = "D"
$ws.Range("K136").Value = 2.8
$ws.Range("L136").Value = 3.3
$ws.Range("M136").Value = 2.3
$ws.Range("N136").Value = 2.55
$ws.Range("O136").Value = 3.2
$ws.Range("P136").Value = 2.6
$ws.Range("Q136").Value = 0
$ws.Range("R136").Value = 1.875
$ws.Range("S136").Value = 1.975
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 2.025
$ws.Range("V136").Value = 1.825
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = 2.2
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = 0
$ws.Range("AA136").Value = -0
$ws.Range("AB136").Value = -1
$ws.Range("AC136").Value = 0.825

# Row 137
$ws.Range("A137").Value = 135
$ws.Range("B137").Value = 6788938
$ws.Range("C137").Value = "Croatia HNL"
$ws.Range("D137").Value = "Croatia HNL"
$ws.Range("E137").Value = 45380.58333333334
$ws.Range("F137").Value = "Slaven Belupo"
$ws.Range("G137").Value = "NK Osijek"
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "A"
$ws.Range("K137").Value = 2.625
$ws.Range("L137").Value = 3.4
$ws.Range("M137").Value = 2.4
$ws.Range("N137").Value = 3.5
$ws.Range("O137").Value = 3.3
$ws.Range("P137").Value = 2
$ws.Range("Q137").Value = 0.5
$ws.Range("R137").Value = 1.8
$ws.Range("S137").Value = 2.05
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.975
$ws.Range("V137").Value = 1.875
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 1
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 1.05
$ws.Range("AB137").Value = -1
$ws.Range("AC137").Value = 0.875

# Row 138
$ws.Range("A138").Value = 136
$ws.Range("B138").Value = 6788936
$ws.Range("C138").Value = "Croatia HNL"
$ws.Range("D138").Value = "Croatia HNL"
$ws.Range("E138").Value = 45381.45138888889
$ws.Range("F138").Value = "NK Varazdin"
$ws.Range("G138").Value = "HNK Gorica"
$ws.Range("H138").Value = 2
$ws.Range("I138").Value = 4
$ws.Range("J138").Value = "A"
$ws.Range("K138").Value = 2.25
$ws.Range("L138").Value = 3.3
$ws.Range("M138").Value = 2.875
$ws.Range("N138").Value = 2.25
$ws.Range("O138").Value = 3.1
$ws.Range("P138").Value = 3.1
$ws.Range("Q138").Value = -0.25
$ws.Range("R138").Value = 1.95
$ws.Range("S138").Value = 1.9
$ws.Range("T138").Value = 2.25
$ws.Range("U138").Value = 2
$ws.Range("V138").Value = 1.85
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = 2.1
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.8999999999999999
$ws.Range("AB138").Value = 1
$ws.Range("AC138").Value = -1

# Row 139
$ws.Range("A139").Value = 137
$ws.Range("B139").Value = 6769306
$ws.Range("C139").Value = "Croatia HNL"
$ws.Range("D139").Value = "Croatia HNL"
$ws.Range("E139").Value = 45381.54166666666
$ws.Range("F139").Value = "NK Rudes"
$ws.Range("G139").Value = "HNK Rijeka"
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 3
$ws.Range("J139").Value = "A"
$ws.Range("K139").Value = 11
$ws.Range("L139").Value = 6
$ws.Range("M139").Value = 1.2
$ws.Range("N139").Value = 21
$ws.Range("O139").Value = 9.5
$ws.Range("P139").Value = 1.083
$ws.Range("Q139").Value = 2.5
$ws.Range("R139").Value = 1.85
$ws.Range("S139").Value = 2
$ws.Range("T139").Value = 3
$ws.Range("U139").Value = 1.825
$ws.Range("V139").Value = 2.025
$ws.Range("W139").Value = -1
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = 0.08299999999999996
$ws.Range("Z139").Value = -1
$ws.Range("AA139").Value = 1
$ws.Range("AB139").Value = 0
$ws.Range("AC139").Value = -0

# Row 140
$ws.Range("A140").Value = 138
$ws.Range("B140").Value = 6788937
$ws.Range("C140").Value = "Croatia HNL"
$ws.Range("D140").Value = "Croatia HNL"
$ws.Range("E140").Value = 45381.64583333334
$ws.Range("F140").Value = "Hajduk Split"
$ws.Range("G140").Value = "Dinamo Zagreb"
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = "A"
$ws.Range("K140").Value = 2.4
$ws.Range("L140").Value = 3.4
$ws.Range("M140").Value = 2.625
$ws.Range("N140").Value = 2.625
$ws.Range("O140").Value = 3.2
$ws.Range("P140").Value = 2.55
$ws.Range("Q140").Value = 0
$ws.Range("R140").Value = 1.95
$ws.Range("S140").Value = 1.9
$ws.Range("T140").Value = 2.25
$ws.Range("U140").Value = 2
$ws.Range("V140").Value = 1.85
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = 1.55
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.8999999999999999
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.8500000000000001
